$wb = $excel.ActiveWorkbook

# Sheet "view_awarding_institutes" -- B5 "Y" -> "N"
$ws16 = $wb.Worksheets.Item("view_awarding_institutes")
$ws16.Range("B5").Value = "N"
$ws16.Range("B5").Select() | Out-Null

# Sheet "edit_awarding_institute" -- D5 "Y" -> "N"
$ws19 = $wb.Worksheets.Item("edit_awarding_institute")
$ws19.Range("D5").Value = "N"
$ws19.Range("D5").Select() | Out-Null

# Sheet "filter_classrooms" -- fix out-of-bound row indexes / keyword
# B2 "CLA" -> "R-005", A4 6->5, A5 7->6, A6 8->7
$ws24 = $wb.Worksheets.Item("filter_classrooms")
$ws24.Range("B2").Value = "R-005"
$ws24.Range("A4").Value = "5"
$ws24.Range("A5").Value = "6"
$ws24.Range("A6").Value = "7"
$ws24.Range("A6").Select() | Out-Null

# Make filter_classrooms the active sheet/tab
$ws24.Activate() | Out-Null
